{"js": "// Localize / refresh the Korean marketing copy for the Mystic Spice Chai\n// product description. Each change below finds the exact existing run\n// text with Body.search (an exact, case-sensitive match) and replaces it\n// in place so surrounding run formatting (fonts, language tags, etc.) is\n// preserved.\n\nasync function replaceOnce(context, oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Title: localize the Korean transliteration of the product name to the\n//    mixed English/Korean brand form used everywhere else in the doc.\nawait replaceOnce(\n  context,\n  \": \ubbf8\uc2a4\ud2f1 \uc2a4\ud30c\uc774\uc2a4 \ud504\ub9ac\ubbf8\uc5c4 \ucc28\uc774 \ucc28\",\n  \": Mystic Spice \ud504\ub9ac\ubbf8\uc5c4 \ucc28\uc774 \ud2f0\"\n);\n\n// 2. \"\uc8fc\uc694 \ud2b9\uc9d5:\" (\"Key Features:\") -> \"\uc8fc\uc694 \uae30\ub2a5:\" (\"Key Functions:\"),\n//    and make that heading run bold (it was incorrectly non-bold, unlike\n//    every other section heading in the document).\nconst keyFeatures = context.document.body.search(\"\uc8fc\uc694 \ud2b9\uc9d5:\", { matchCase: true });\nkeyFeatures.load(\"items\");\nawait context.sync();\nif (keyFeatures.items.length === 0) {\n  throw new Error(\"Text not found: \uc8fc\uc694 \ud2b9\uc9d5:\");\n}\nconst keyFeaturesRange = keyFeatures.items[0];\nkeyFeaturesRange.insertText(\"\uc8fc\uc694 \uae30\ub2a5:\", Word.InsertLocation.replace);\nawait context.sync();\n// Re-search for the freshly inserted text so we target the right range,\n// then flip it to bold.\nconst keyFunctions = context.document.body.search(\"\uc8fc\uc694 \uae30\ub2a5:\", { matchCase: true });\nkeyFunctions.load(\"items\");\nawait context.sync();\nkeyFunctions.items[0].font.bold = true;\nawait context.sync();\n\n// 3. \"Authentic Blend\" bullet detail.\nawait replaceOnce(\n  context,\n  \": \uc800\ud76c \ucc28\uc774\ub294 \ud504\ub9ac\ubbf8\uc5c4 \ud64d\ucc28 \uc78e\uacfc \uacc4\ud53c, \uce74\ub2e4\ubaac, \uc815\ud5a5, \uc0dd\uac15, \ud6c4\ucd94 \ub4f1 \ub2e4\uc591\ud55c \uc9c0\uc0c1 \ud5a5\uc2e0\ub8cc\uc758 \uc870\ud654\ub85c\uc6b4 \uc870\ud569\uc785\ub2c8\ub2e4. \",\n  \": \uc800\ud76c \ucc28\uc774\ub294 \ud504\ub9ac\ubbf8\uc5c4 \ud64d\ucc28 \uc78e\uacfc \uacc4\ud53c, \uce74\ub2e4\ubaac, \uc815\ud5a5, \uc0dd\uac15, \ud6c4\ucd94 \ub4f1 \ub2e4\uc591\ud55c \uc2dc\uadf8\ub2c8\ucc98 \uac00\ub8e8 \ud5a5\uc2e0\ub8cc\uc640\uc758 \uc870\ud654\ub85c\uc6b4 \ubbf9\uc2a4\ub85c \uc774\ub8e8\uc5b4\uc9d1\ub2c8\ub2e4. \"\n);\n\n// 4. \"Health-Boosting Ingredients\" bullet detail.\nawait replaceOnce(\n  context,\n  \": \uc2e0\ube44 \ud5a5\uc2e0\ub8cc \ucc28\uc774 \ucc28\uc758 \uac01 \uc131\ubd84\uc740 \ucc9c\uc5f0 \uac74\uac15\uc0c1\uc758 \uc774\uc810\uc744 \uc704\ud574 \uc120\ud0dd\ub429\ub2c8\ub2e4. \",\n  \": Mystic Spice \ucc28\uc774 \ud2f0\uc758 \uac01 \uc131\ubd84\uc740 \ucc9c\uc5f0\uc758 \uac74\uac15 \ud61c\ud0dd\uc5d0 \uae30\ubc18\ud558\uc5ec \uc5c4\uc120\ub429\ub2c8\ub2e4. \"\n);\n\n// 5. \"Rich Aroma and Flavor\" bullet detail.\nawait replaceOnce(\n  context,\n  \": \ub530\ub73b\ud558\uace0 \ub9e4\uc6b4 \ud5a5\uae30\uc640 \uc6b0\ub9ac\uc758 \ucc28\uc774\uc758 \uae4a\uace0 \uc0c1\ucf8c\ud55c \ub9db\uc740 \ud558\ub8e8\ub97c \uc2dc\uc791\ud558\uac70\ub098 \uc800\ub141\uc5d0 \uae34\uc7a5\uc744 \ud480 \uc218\uc788\ub294 \uc644\ubcbd\ud55c \uc74c\ub8cc\uc785\ub2c8\ub2e4. \",\n  \": \uc800\ud76c \ucc28\uc774\uc758 \ub530\ub73b\ud558\uace0 \ub9e4\uc6b4 \ud5a5\uae30\uc640 \uae4a\uace0 \uc0c1\ucf8c\ud55c \ub9db\uc740 \ud558\ub8e8\ub97c \uc2dc\uc791\ud558\uac70\ub098 \uc800\ub141\uc5d0 \uae34\uc7a5\uc744 \ud480\uae30\uc5d0 \uc644\ubcbd\ud55c \uc74c\ub8cc\uc758 \uc870\uac74\uc785\ub2c8\ub2e4. \"\n);\n\n// 6. \"Versatile Brewing Options\" bullet heading.\nawait replaceOnce(\n  context,\n  \"\ub2e4\uc7ac\ub2e4\ub2a5\ud55c \uc591\uc870 \uc635\uc158\",\n  \"\ub2e4\uc591\ud55c \ube0c\ub8e8\uc789 \uc635\uc158\"\n);\n\n// 7. \"Versatile Brewing Options\" bullet detail.\nawait replaceOnce(\n  context,\n  \": \ucc28\uc774 \uae40\uc774 \ub728\uac70\uc6cc\uc9c0\uac70\ub098, \uc0c1\ucf8c\ud55c \uc544\uc774\uc2a4 \ud2f0\ub85c, \ud06c\ub9ac\ubbf8\ud55c \ub77c\ub5bc\ub97c \uc88b\uc544\ud558\ub4e0, \uc800\ud76c \ube14\ub80c\ub4dc\ub294 \ubaa8\ub4e0 \ucde8\ud5a5\uc5d0 \ub9de\uac8c \ub2e4\uc7ac\ub2e4\ub2a5\ud569\ub2c8\ub2e4. \",\n  \": \uae40\uc774 \ub728\uac81\uac8c \uc62c\ub77c\uc624\ub294 \ucc28\uc774\ub098, \uc0c1\ucf8c\ud55c \uc544\uc774\uc2a4 \ud2f0, \ud06c\ub9ac\ubbf8\ud55c \ub77c\ub5bc \ub4f1, \uc800\ud76c \ube14\ub80c\ub4dc\ub294 \ubaa8\ub4e0 \ucde8\ud5a5\uc5d0 \ub9de\ub294 \ub2e4\uc591\uc131\uc744 \uc81c\uacf5\ud569\ub2c8\ub2e4. \"\n);\n\n// 8. \"Sustainable Sourcing\" bullet heading.\nawait replaceOnce(\n  context,\n  \"\uc9c0\uc18d \uac00\ub2a5\ud55c \uacf5\uae09\",\n  \"\uc9c0\uc18d \uac00\ub2a5\ud55c \uc18c\uc2f1\"\n);\n\n// 9. \"Sustainable Sourcing\" bullet detail.\nawait replaceOnce(\n  context,\n  \": \uc9c0\uc18d \uac00\ub2a5\uc131\uc744 \uc704\ud574 \ucd5c\uc120\uc744 \ub2e4\ud558\uace0 \uc788\uc73c\uba70, \uc6b0\ub9ac\ub294 \uc720\uae30\ub18d \ub18d\uc5c5\uc744 \uc2e4\ucc9c\ud558\ub294 \uc18c\uaddc\ubaa8 \ub18d\uc7a5\uc5d0\uc11c \uc7ac\ub8cc\ub97c \uacf5\uae09\ud558\uc5ec \ucd5c\uace0\uc758 \ud488\uc9c8\ubfd0\ub9cc \uc544\ub2c8\ub77c \uc9c0\uad6c\uc758 \ubcf5\uc9c0\ub97c \ubcf4\uc7a5\ud569\ub2c8\ub2e4.\",\n  \": \uc800\ud76c\ub294 \uc9c0\uc18d \uac00\ub2a5\uc131\uc744 \uc704\ud574 \ucd5c\uc120\uc744 \ub2e4\ud558\uace0 \uc788\uc73c\uba70, \uc720\uae30\ub18d \ub18d\uc5c5\uc744 \uc2e4\ucc9c\ud558\ub294 \uc18c\uaddc\ubaa8 \ub18d\uc7a5\uc5d0\uc11c \uc7ac\ub8cc\ub97c \uacf5\uae09\ubc1b\uc74c\uc73c\ub85c\uc368 \ucd5c\uace0\uc758 \ud488\uc9c8\ubfd0\ub9cc \uc544\ub2c8\ub77c \uc9c0\uad6c\uc758 \uac74\uac15\uc5d0\ub3c4 \uc774\ubc14\uc9c0\ud558\uace0 \uc788\uc2b5\ub2c8\ub2e4.\"\n);\n\n// 10. \"Elegant Packaging\" bullet detail.\nawait replaceOnce(\n  context,\n  \": \uc2e0\ube44\uc8fc\uc758 \ud5a5\uc2e0\ub8cc \ucc28\uc774 \ucc28\ub294 \uc544\ub984\ub2f5\uac8c \ub514\uc790\uc778\ub41c \uce5c\ud658\uacbd \ud3ec\uc7a5\uc7ac\ub85c \uc81c\uacf5\ub418\uba70, \ucc28 \uc560\ud638\uac00\ub4e4\uc5d0\uac8c \uc774\uc0c1\uc801\uc778 \uc120\ubb3c\uc774\uac70\ub098 \ud638\ud654\ub85c\uc6b4 \uac04\uc2dd\uc785\ub2c8\ub2e4.\",\n  \": Mystic Spice \ucc28\uc774 \ud2f0\ub294 \uc544\ub984\ub2f5\uac8c \ub514\uc790\uc778\ub41c \uce5c\ud658\uacbd \ud3ec\uc7a5\uc7ac\ub85c \uc81c\uacf5\ub418\ubbc0\ub85c \ucc28 \uc560\ud638\uac00\ub4e4\uc5d0\uac8c \ubfd0\ub9cc \uc544\ub2c8\ub77c \uc790\uc2e0\uc5d0\uac8c\ub3c4 \uc774\uc0c1\uc801\uc774\uace0 \uace0\uae09\uc2a4\ub7ec\uc6b4 \uc120\ubb3c\uc774 \ub420 \uc218 \uc788\uc2b5\ub2c8\ub2e4.\"\n);\n\n// 11. \"Customer Satisfaction Guarantee\" bullet detail.\nawait replaceOnce(\n  context,\n  \": Microsoft\ub294 \uc81c\ud488 \ub4a4\uc5d0 \uc11c\uc11c \ub9cc\uc871\ub3c4 \ubcf4\uc7a5\uc744 \uc81c\uacf5\ud569\ub2c8\ub2e4. \",\n  \": \uc800\ud76c\ub294 \uc81c\ud488\uc744 \ub4b7\ubc1b\uce68\ud558\uba70 \uace0\uac1d \ub9cc\uc871\uc744 \ubcf4\uc7a5\ud569\ub2c8\ub2e4. \"\n);\n\n// 12. \"Ideal For\" bullet detail.\nawait replaceOnce(\n  context,\n  \": \ucc28 \uc560\ud638\uac00, \uac74\uac15\uc5d0 \ubbfc\uac10\ud55c \uac1c\uc778, \ub530\ub73b\ud558\uace0 \ub9e4\uc6b4 \uc74c\ub8cc \uc560\ud638\uac00, \uc804\ud1b5\uc801\uc778 \uc778\ub3c4 \ucc28\uc774\uc758 \ud48d\ubd80\ud55c \ub9db\uc744 \ud0d0\uad6c\ud558\uace0\uc790\ud558\ub294 \uc0ac\ub78c.\",\n  \": \ucc28 \uc560\ud638\uac00, \uac74\uac15\uc5d0 \ubbfc\uac10\ud55c \uc0ac\ub78c, \ub530\ub73b\ud558\uace0 \ub9e4\uc6b4 \uc74c\ub8cc \uc560\ud638\uac00, \uc804\ud1b5\uc801\uc778 \uc778\ub3c4 \ucc28\uc774\uc758 \ud48d\ubd80\ud55c \ub9db\uc744 \ud0d0\uad6c\ud558\uace0\uc790\ud558\ub294 \ubaa8\ub4e0 \uc0ac\ub78c.\"\n);\n", "ps1": "# Localize / refresh the Korean marketing copy for the Mystic Spice Chai\n# product description. Each Find/Replace below targets the exact existing\n# run text so the surrounding run formatting (fonts, language tags, etc.)\n# is preserved; wdReplaceOne (2) limits each call to the single match.\n\n$d = $word.ActiveDocument\n\n# 1. Title: localize the Korean transliteration of the product name to the\n#    mixed English/Korean brand form used everywhere else in the doc.\n$d.Content.Find.Execute(\": \ubbf8\uc2a4\ud2f1 \uc2a4\ud30c\uc774\uc2a4 \ud504\ub9ac\ubbf8\uc5c4 \ucc28\uc774 \ucc28\", $false, $false, $false, $false, $false, $true, 1, $false, \": Mystic Spice \ud504\ub9ac\ubbf8\uc5c4 \ucc28\uc774 \ud2f0\", 2)\n\n# 2. \"\uc8fc\uc694 \ud2b9\uc9d5:\" (\"Key Features:\") -> \"\uc8fc\uc694 \uae30\ub2a5:\" (\"Key Functions:\").\n$d.Content.Find.Execute(\"\uc8fc\uc694 \ud2b9\uc9d5:\", $false, $false, $false, $false, $false, $true, 1, $false, \"\uc8fc\uc694 \uae30\ub2a5:\", 2)\n\n# Make that heading run bold (it was incorrectly non-bold, unlike every\n# other section heading in the document).\n$headingRange = $d.Content\n$headingRange.Find.Execute(\"\uc8fc\uc694 \uae30\ub2a5:\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n$headingRange.Font.Bold = 1\n\n# 3. \"Authentic Blend\" bullet detail.\n$d.Content.Find.Execute(\": \uc800\ud76c \ucc28\uc774\ub294 \ud504\ub9ac\ubbf8\uc5c4 \ud64d\ucc28 \uc78e\uacfc \uacc4\ud53c, \uce74\ub2e4\ubaac, \uc815\ud5a5, \uc0dd\uac15, \ud6c4\ucd94 \ub4f1 \ub2e4\uc591\ud55c \uc9c0\uc0c1 \ud5a5\uc2e0\ub8cc\uc758 \uc870\ud654\ub85c\uc6b4 \uc870\ud569\uc785\ub2c8\ub2e4. \", $false, $false, $false, $false, $false, $true, 1, $false, \": \uc800\ud76c \ucc28\uc774\ub294 \ud504\ub9ac\ubbf8\uc5c4 \ud64d\ucc28 \uc78e\uacfc \uacc4\ud53c, \uce74\ub2e4\ubaac, \uc815\ud5a5, \uc0dd\uac15, \ud6c4\ucd94 \ub4f1 \ub2e4\uc591\ud55c \uc2dc\uadf8\ub2c8\ucc98 \uac00\ub8e8 \ud5a5\uc2e0\ub8cc\uc640\uc758 \uc870\ud654\ub85c\uc6b4 \ubbf9\uc2a4\ub85c \uc774\ub8e8\uc5b4\uc9d1\ub2c8\ub2e4. \", 2)\n\n# 4. \"Health-Boosting Ingredients\" bullet detail.\n$d.Content.Find.Execute(\": \uc2e0\ube44 \ud5a5\uc2e0\ub8cc \ucc28\uc774 \ucc28\uc758 \uac01 \uc131\ubd84\uc740 \ucc9c\uc5f0 \uac74\uac15\uc0c1\uc758 \uc774\uc810\uc744 \uc704\ud574 \uc120\ud0dd\ub429\ub2c8\ub2e4. \", $false, $false, $false, $false, $false, $true, 1, $false, \": Mystic Spice \ucc28\uc774 \ud2f0\uc758 \uac01 \uc131\ubd84\uc740 \ucc9c\uc5f0\uc758 \uac74\uac15 \ud61c\ud0dd\uc5d0 \uae30\ubc18\ud558\uc5ec \uc5c4\uc120\ub429\ub2c8\ub2e4. \", 2)\n\n# 5. \"Rich Aroma and Flavor\" bullet detail.\n$d.Content.Find.Execute(\": \ub530\ub73b\ud558\uace0 \ub9e4\uc6b4 \ud5a5\uae30\uc640 \uc6b0\ub9ac\uc758 \ucc28\uc774\uc758 \uae4a\uace0 \uc0c1\ucf8c\ud55c \ub9db\uc740 \ud558\ub8e8\ub97c \uc2dc\uc791\ud558\uac70\ub098 \uc800\ub141\uc5d0 \uae34\uc7a5\uc744 \ud480 \uc218\uc788\ub294 \uc644\ubcbd\ud55c \uc74c\ub8cc\uc785\ub2c8\ub2e4. \", $false, $false, $false, $false, $false, $true, 1, $false, \": \uc800\ud76c \ucc28\uc774\uc758 \ub530\ub73b\ud558\uace0 \ub9e4\uc6b4 \ud5a5\uae30\uc640 \uae4a\uace0 \uc0c1\ucf8c\ud55c \ub9db\uc740 \ud558\ub8e8\ub97c \uc2dc\uc791\ud558\uac70\ub098 \uc800\ub141\uc5d0 \uae34\uc7a5\uc744 \ud480\uae30\uc5d0 \uc644\ubcbd\ud55c \uc74c\ub8cc\uc758 \uc870\uac74\uc785\ub2c8\ub2e4. \", 2)\n\n# 6. \"Versatile Brewing Options\" bullet heading.\n$d.Content.Find.Execute(\"\ub2e4\uc7ac\ub2e4\ub2a5\ud55c \uc591\uc870 \uc635\uc158\", $false, $false, $false, $false, $false, $true, 1, $false, \"\ub2e4\uc591\ud55c \ube0c\ub8e8\uc789 \uc635\uc158\", 2)\n\n# 7. \"Versatile Brewing Options\" bullet detail.\n$d.Content.Find.Execute(\": \ucc28\uc774 \uae40\uc774 \ub728\uac70\uc6cc\uc9c0\uac70\ub098, \uc0c1\ucf8c\ud55c \uc544\uc774\uc2a4 \ud2f0\ub85c, \ud06c\ub9ac\ubbf8\ud55c \ub77c\ub5bc\ub97c \uc88b\uc544\ud558\ub4e0, \uc800\ud76c \ube14\ub80c\ub4dc\ub294 \ubaa8\ub4e0 \ucde8\ud5a5\uc5d0 \ub9de\uac8c \ub2e4\uc7ac\ub2e4\ub2a5\ud569\ub2c8\ub2e4. \", $false, $false, $false, $false, $false, $true, 1, $false, \": \uae40\uc774 \ub728\uac81\uac8c \uc62c\ub77c\uc624\ub294 \ucc28\uc774\ub098, \uc0c1\ucf8c\ud55c \uc544\uc774\uc2a4 \ud2f0, \ud06c\ub9ac\ubbf8\ud55c \ub77c\ub5bc \ub4f1, \uc800\ud76c \ube14\ub80c\ub4dc\ub294 \ubaa8\ub4e0 \ucde8\ud5a5\uc5d0 \ub9de\ub294 \ub2e4\uc591\uc131\uc744 \uc81c\uacf5\ud569\ub2c8\ub2e4. \", 2)\n\n# 8. \"Sustainable Sourcing\" bullet heading.\n$d.Content.Find.Execute(\"\uc9c0\uc18d \uac00\ub2a5\ud55c \uacf5\uae09\", $false, $false, $false, $false, $false, $true, 1, $false, \"\uc9c0\uc18d \uac00\ub2a5\ud55c \uc18c\uc2f1\", 2)\n\n# 9. \"Sustainable Sourcing\" bullet detail.\n$d.Content.Find.Execute(\": \uc9c0\uc18d \uac00\ub2a5\uc131\uc744 \uc704\ud574 \ucd5c\uc120\uc744 \ub2e4\ud558\uace0 \uc788\uc73c\uba70, \uc6b0\ub9ac\ub294 \uc720\uae30\ub18d \ub18d\uc5c5\uc744 \uc2e4\ucc9c\ud558\ub294 \uc18c\uaddc\ubaa8 \ub18d\uc7a5\uc5d0\uc11c \uc7ac\ub8cc\ub97c \uacf5\uae09\ud558\uc5ec \ucd5c\uace0\uc758 \ud488\uc9c8\ubfd0\ub9cc \uc544\ub2c8\ub77c \uc9c0\uad6c\uc758 \ubcf5\uc9c0\ub97c \ubcf4\uc7a5\ud569\ub2c8\ub2e4.\", $false, $false, $false, $false, $false, $true, 1, $false, \": \uc800\ud76c\ub294 \uc9c0\uc18d \uac00\ub2a5\uc131\uc744 \uc704\ud574 \ucd5c\uc120\uc744 \ub2e4\ud558\uace0 \uc788\uc73c\uba70, \uc720\uae30\ub18d \ub18d\uc5c5\uc744 \uc2e4\ucc9c\ud558\ub294 \uc18c\uaddc\ubaa8 \ub18d\uc7a5\uc5d0\uc11c \uc7ac\ub8cc\ub97c \uacf5\uae09\ubc1b\uc74c\uc73c\ub85c\uc368 \ucd5c\uace0\uc758 \ud488\uc9c8\ubfd0\ub9cc \uc544\ub2c8\ub77c \uc9c0\uad6c\uc758 \uac74\uac15\uc5d0\ub3c4 \uc774\ubc14\uc9c0\ud558\uace0 \uc788\uc2b5\ub2c8\ub2e4.\", 2)\n\n# 10. \"Elegant Packaging\" bullet detail.\n$d.Content.Find.Execute(\": \uc2e0\ube44\uc8fc\uc758 \ud5a5\uc2e0\ub8cc \ucc28\uc774 \ucc28\ub294 \uc544\ub984\ub2f5\uac8c \ub514\uc790\uc778\ub41c \uce5c\ud658\uacbd \ud3ec\uc7a5\uc7ac\ub85c \uc81c\uacf5\ub418\uba70, \ucc28 \uc560\ud638\uac00\ub4e4\uc5d0\uac8c \uc774\uc0c1\uc801\uc778 \uc120\ubb3c\uc774\uac70\ub098 \ud638\ud654\ub85c\uc6b4 \uac04\uc2dd\uc785\ub2c8\ub2e4.\", $false, $false, $false, $false, $false, $true, 1, $false, \": Mystic Spice \ucc28\uc774 \ud2f0\ub294 \uc544\ub984\ub2f5\uac8c \ub514\uc790\uc778\ub41c \uce5c\ud658\uacbd \ud3ec\uc7a5\uc7ac\ub85c \uc81c\uacf5\ub418\ubbc0\ub85c \ucc28 \uc560\ud638\uac00\ub4e4\uc5d0\uac8c \ubfd0\ub9cc \uc544\ub2c8\ub77c \uc790\uc2e0\uc5d0\uac8c\ub3c4 \uc774\uc0c1\uc801\uc774\uace0 \uace0\uae09\uc2a4\ub7ec\uc6b4 \uc120\ubb3c\uc774 \ub420 \uc218 \uc788\uc2b5\ub2c8\ub2e4.\", 2)\n\n# 11. \"Customer Satisfaction Guarantee\" bullet detail.\n$d.Content.Find.Execute(\": Microsoft\ub294 \uc81c\ud488 \ub4a4\uc5d0 \uc11c\uc11c \ub9cc\uc871\ub3c4 \ubcf4\uc7a5\uc744 \uc81c\uacf5\ud569\ub2c8\ub2e4. \", $false, $false, $false, $false, $false, $true, 1, $false, \": \uc800\ud76c\ub294 \uc81c\ud488\uc744 \ub4b7\ubc1b\uce68\ud558\uba70 \uace0\uac1d \ub9cc\uc871\uc744 \ubcf4\uc7a5\ud569\ub2c8\ub2e4. \", 2)\n\n# 12. \"Ideal For\" bullet detail.\n$d.Content.Find.Execute(\": \ucc28 \uc560\ud638\uac00, \uac74\uac15\uc5d0 \ubbfc\uac10\ud55c \uac1c\uc778, \ub530\ub73b\ud558\uace0 \ub9e4\uc6b4 \uc74c\ub8cc \uc560\ud638\uac00, \uc804\ud1b5\uc801\uc778 \uc778\ub3c4 \ucc28\uc774\uc758 \ud48d\ubd80\ud55c \ub9db\uc744 \ud0d0\uad6c\ud558\uace0\uc790\ud558\ub294 \uc0ac\ub78c.\", $false, $false, $false, $false, $false, $true, 1, $false, \": \ucc28 \uc560\ud638\uac00, \uac74\uac15\uc5d0 \ubbfc\uac10\ud55c \uc0ac\ub78c, \ub530\ub73b\ud558\uace0 \ub9e4\uc6b4 \uc74c\ub8cc \uc560\ud638\uac00, \uc804\ud1b5\uc801\uc778 \uc778\ub3c4 \ucc28\uc774\uc758 \ud48d\ubd80\ud55c \ub9db\uc744 \ud0d0\uad6c\ud558\uace0\uc790\ud558\ub294 \ubaa8\ub4e0 \uc0ac\ub78c.\", 2)\n"}
